$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the README sheet: insert a new row describing the new
#    OUTPUT_SHEETS tab, just after the STANDARDS description line.
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("README")
$readme.Rows.Item(12).Insert()
$readme.Range("B12").Value = "OUTPUT_SHEETS: A description of the data contained in the exported data"

# ---------------------------------------------------------------------------
# 2. Add a new OUTPUT_SHEETS worksheet at the end of the workbook (after
#    STANDARDS) describing the sheets contained in the tool's export files.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "OUTPUT_SHEETS"

$ws.Range("A1").Value = "Sheet name"
$ws.Range("B1").Value = "Description"

$ws.Range("A2").Value = "PROBE_DATA"
$ws.Range("B2").Value = "Original data gathered from the automated probes."

$ws.Range("A3").Value = "MANUAL_FIELD"
$ws.Range("B3").Value = "Original data gathered manually from the field."

$ws.Range("A4").Value = "LAB_DATA"
$ws.Range("B4").Value = "Original data from the labs."

$ws.Range("A5").Value = "EA_DATA"
$ws.Range("B5").Value = "Original data downloaded from WIMS."

$ws.Range("A6").Value = "OUTWITH_LOD"
$ws.Range("B6").Value = 'Results identified as being outwith the limit of detection (based on the presence of ">" or "<" in the results field).'

$ws.Range("A7").Value = "EXCLUSIONS"
$ws.Range("B7").Value = "Results that have been excluded from analysis as a results of being outwith the natural range of the parameter."

$ws.Range("A8").Value = "OUTLIERS"
$ws.Range("B8").Value = "Results that have been identified as potential outliers."

$ws.Range("A9").Value = "SUMMARY"
$ws.Range("B9").Value = "A statistical summary of results by site and parameter."

$ws.Range("A10").Value = "FULL_DATA"
$ws.Range("B10").Value = "Full processed results from the analysis."

$ws.Range("A12").Value = "Notes:"
$ws.Range("B13").Value = "Sheets PROBE_DATA, MANUAL_FIELD, LAB_DATA and EA_DATA will be blank if data from these sources was not supplied."

[void]$ws.Range("A1:B13").Select()

# ---------------------------------------------------------------------------
# 3. Restore the original active sheet / selection so the workbook still
#    opens on the README tab.
# ---------------------------------------------------------------------------
[void]$readme.Activate()
[void]$readme.Range("A1").Select()
